$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user selected entire rows 19:20 and deleted them (shift cells up),
# which moves the "Jenis mimpi" / dream-related block up by two rows.
$ws.Rows("19:20").Delete()

# Excel re-saves the file and, in doing so, de-duplicates the cell style
# table: several style records used only by the Qur'an-reference cells in
# column G:J were byte-for-byte identical (same font, same "no border"
# formatting) and collapse into a single shared style. Re-apply the format
# already used by G5:J5 to those ranges so the redundant styles merge away.
$ws.Range("G5:J5").Copy()
$ws.Range("G6:J6").PasteSpecial(-4122)
$ws.Range("G8:J9").PasteSpecial(-4122)
$ws.Range("G18:J18").PasteSpecial(-4122)
$ws.Range("G21:J21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reflect the selection the user was left with after deleting rows 19:20.
$ws.Range("A19:XFD20").Select()

$wb.Save()
